# "updated clean stat files"
# Sets column B (playerid) to a constant 32 for every game row (rows 3-83 on
# Sheet1 — row 2 is already 32 and is left untouched), including filling in
# the previously-empty B cells for rows 29, 59 and 64-81. Finally restores
# the last-used cell selection to D8, matching the saved sheetView state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3:B83").Value = 32

$ws.Range("D8").Select()
